# New crime data collected — refresh the 33rd Precinct weekly CompStat sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header text: bump the bulletin "Volume/Number" and the reporting week.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  35"
$ws.Range("C9").Value = "Report Covering the Week  8/28/2023  Through  9/3/2023"

# ---------------------------------------------------------------------------
# 2) Stable "donor" cells used purely to copy a cell's number format (style)
#    onto a cell whose underlying type (text placeholder <-> number) changes.
#    These donors live outside the edited row block (14-29) and are never
#    themselves modified by this script, so they stay valid for the whole run.
# ---------------------------------------------------------------------------
$donorTextZero   = $ws.Range("G14")   # style 14, text "0"
$donorTextStar   = $ws.Range("H14")   # style 14, text "***.*"
$donorNumberPlain = $ws.Range("E36")  # style 15, plain integer number format
$donorNumberPct   = $ws.Range("K36")  # style 16, percent-change number format

function Set-TextZero($rng) {
    $rng.Value = "'0"
    $donorTextZero.Copy()
    $rng.PasteSpecial(-4122)
}

function Set-TextStar($rng) {
    $rng.Value = "'***.*"
    $donorTextStar.Copy()
    $rng.PasteSpecial(-4122)
}

function Set-NumberPlain($rng, $val) {
    $donorNumberPlain.Copy()
    $rng.PasteSpecial(-4122)
    $rng.Value = $val
}

function Set-NumberPct($rng, $val) {
    $donorNumberPct.Copy()
    $rng.PasteSpecial(-4122)
    $rng.Value = $val
}

# ---------------------------------------------------------------------------
# 3) Cells that flip from a number to the "0"/"***.*" text placeholder.
# ---------------------------------------------------------------------------
Set-TextZero $ws.Range("C14")
Set-TextZero $ws.Range("D15")
Set-TextStar $ws.Range("E15")
Set-TextZero $ws.Range("C20")
Set-TextZero $ws.Range("C22")
Set-TextZero $ws.Range("D22")
Set-TextStar $ws.Range("E22")
Set-TextZero $ws.Range("D28")
Set-TextStar $ws.Range("E28")
Set-TextZero $ws.Range("F28")
Set-TextZero $ws.Range("D29")
Set-TextStar $ws.Range("E29")
Set-TextZero $ws.Range("F29")

# ---------------------------------------------------------------------------
# 4) Cells that flip from the "0"/"***.*" text placeholder to a real number.
# ---------------------------------------------------------------------------
Set-NumberPlain $ws.Range("D20") 7
Set-NumberPct   $ws.Range("E20") -100
Set-NumberPlain $ws.Range("C23") 1
Set-NumberPlain $ws.Range("C27") 1
Set-NumberPlain $ws.Range("D27") 1
Set-NumberPct   $ws.Range("E27") 0

# ---------------------------------------------------------------------------
# 5) Plain value refreshes (style / type unchanged) — weekly, 28-day, YTD and
#    2-year crime-complaint counts/percentages for rows 14-29.
# ---------------------------------------------------------------------------
$ws.Range("F14").Value = 1
$ws.Range("L14").Value = 66.666666666666

$ws.Range("G15").Value = 2

$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 13.333333333333
$ws.Range("I16").Value = 112
$ws.Range("J16").Value = 132
$ws.Range("K16").Value = -15.151515151515
$ws.Range("L16").Value = -29.113924050632
$ws.Range("M16").Value = -27.741935483871

$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 40
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 14.285714285714
$ws.Range("I17").Value = 173
$ws.Range("J17").Value = 198
$ws.Range("K17").Value = -12.626262626262
$ws.Range("L17").Value = -15.609756097561
$ws.Range("M17").Value = 39.516129032258

$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 16.666666666666
$ws.Range("I18").Value = 80
$ws.Range("J18").Value = 66
$ws.Range("K18").Value = 21.212121212121
$ws.Range("L18").Value = -17.525773195876
$ws.Range("M18").Value = 21.212121212121

$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 40
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 19
$ws.Range("H19").Value = 100
$ws.Range("I19").Value = 244
$ws.Range("J19").Value = 207
$ws.Range("K19").Value = 17.874396135265
$ws.Range("L19").Value = 15.094339622641
$ws.Range("M19").Value = 37.078651685393

$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -11.111111111111
$ws.Range("J20").Value = 101
$ws.Range("K20").Value = 0
$ws.Range("M20").Value = 119.565217391304

$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 5
$ws.Range("F21").Value = 95
$ws.Range("G21").Value = 72
$ws.Range("H21").Value = 31.944444444444
$ws.Range("I21").Value = 721
$ws.Range("J21").Value = 717
$ws.Range("K21").Value = 0.557880055788
$ws.Range("L21").Value = -0.961538461538
$ws.Range("M21").Value = 23.883161512027

$ws.Range("F22").Value = 5
$ws.Range("H22").Value = 150
$ws.Range("L22").Value = 18.181818181818
$ws.Range("M22").Value = 30

$ws.Range("I23").Value = 10
$ws.Range("K23").Value = -33.333333333333
$ws.Range("L23").Value = -37.5
$ws.Range("M23").Value = 42.857142857142

$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = -40.625
$ws.Range("F24").Value = 83
$ws.Range("G24").Value = 96
$ws.Range("H24").Value = -13.541666666666
$ws.Range("I24").Value = 659
$ws.Range("J24").Value = 987
$ws.Range("K24").Value = -33.232016210739
$ws.Range("L24").Value = -35.518590998043
$ws.Range("M24").Value = 71.614583333333

$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = 22.222222222222
$ws.Range("I25").Value = 274
$ws.Range("J25").Value = 279
$ws.Range("K25").Value = -1.792114695340
$ws.Range("L25").Value = 1.107011070110
$ws.Range("M25").Value = -17.469879518072

$ws.Range("F26").Value = 2
$ws.Range("H26").Value = -33.333333333333
$ws.Range("I26").Value = 11
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = -38.888888888888
$ws.Range("L26").Value = -35.294117647058

$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 66.666666666666
$ws.Range("I27").Value = 27
$ws.Range("J27").Value = 37
$ws.Range("K27").Value = -27.027027027027
$ws.Range("L27").Value = -15.625

$ws.Range("H28").Value = -100
$ws.Range("L28").Value = -30.769230769230

$ws.Range("H29").Value = -100
$ws.Range("L29").Value = -45.454545454545
